# Update scraped_at timestamps (column K) on the "snapshot" sheet,
# and remove the single data row (row 2) from the "new_injured" sheet
# — a record that was previously marked INJURED_NEW and has since been
# resolved/removed from this run's output.

$wb = $excel.ActiveWorkbook

# --- Sheet "snapshot": refresh column K ("scraped_at") timestamps ---
$snapshot = $wb.Worksheets.Item("snapshot")

$timestamps = @{
    "K2"  = "2025-11-03T15:40:13.087744+00:00";
    "K3"  = "2025-11-03T15:40:15.194452+00:00";
    "K4"  = "2025-11-03T15:40:15.194471+00:00";
    "K5"  = "2025-11-03T15:40:15.194479+00:00";
    "K6"  = "2025-11-03T15:40:15.194487+00:00";
    "K7"  = "2025-11-03T15:40:17.351523+00:00";
    "K8"  = "2025-11-03T15:40:17.351578+00:00";
    "K9"  = "2025-11-03T15:40:17.351604+00:00";
    "K10" = "2025-11-03T15:40:19.611904+00:00";
    "K11" = "2025-11-03T15:40:22.490728+00:00";
    "K12" = "2025-11-03T15:40:22.490759+00:00";
    "K13" = "2025-11-03T15:40:22.490779+00:00";
    "K14" = "2025-11-03T15:40:24.847730+00:00";
    "K15" = "2025-11-03T15:40:24.847758+00:00";
    "K16" = "2025-11-03T15:40:24.847776+00:00";
    "K17" = "2025-11-03T15:40:32.598800+00:00";
    "K18" = "2025-11-03T15:40:34.937654+00:00";
    "K19" = "2025-11-03T15:40:37.424167+00:00";
    "K20" = "2025-11-03T15:40:39.850285+00:00";
    "K21" = "2025-11-03T15:40:39.850317+00:00";
    "K22" = "2025-11-03T15:40:39.850335+00:00";
    "K23" = "2025-11-03T15:40:42.368164+00:00";
    "K24" = "2025-11-03T15:40:42.368198+00:00";
    "K25" = "2025-11-03T15:40:42.368219+00:00";
    "K26" = "2025-11-03T15:40:42.368240+00:00";
    "K27" = "2025-11-03T15:40:44.738644+00:00";
    "K28" = "2025-11-03T15:40:49.938718+00:00";
    "K29" = "2025-11-03T15:40:49.938752+00:00";
    "K30" = "2025-11-03T15:40:49.938770+00:00";
    "K31" = "2025-11-03T15:40:49.938786+00:00";
    "K32" = "2025-11-03T15:40:52.319933+00:00";
    "K33" = "2025-11-03T15:40:52.319965+00:00";
    "K34" = "2025-11-03T15:40:52.319985+00:00";
    "K35" = "2025-11-03T15:40:54.647538+00:00";
    "K36" = "2025-11-03T15:40:54.647570+00:00";
    "K37" = "2025-11-03T15:40:54.647589+00:00";
    "K38" = "2025-11-03T15:40:54.647606+00:00";
    "K39" = "2025-11-03T15:40:54.647623+00:00";
    "K40" = "2025-11-03T15:40:54.647669+00:00";
    "K41" = "2025-11-03T15:40:54.647688+00:00";
    "K42" = "2025-11-03T15:40:54.647708+00:00";
    "K43" = "2025-11-03T15:40:57.622509+00:00";
    "K44" = "2025-11-03T15:40:57.622539+00:00";
    "K45" = "2025-11-03T15:41:02.919375+00:00";
    "K46" = "2025-11-03T15:41:05.312080+00:00";
    "K47" = "2025-11-03T15:41:05.312108+00:00";
    "K48" = "2025-11-03T15:41:05.312125+00:00";
}

foreach ($addr in $timestamps.Keys) {
    $snapshot.Range($addr).Value = $timestamps[$addr]
}

# --- Sheet "new_injured": drop the resolved row (row 2) ---
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows.Item(2).Delete()
